# Auto-generated Excel COM-interop script to apply LDLC price-history diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at FE, shifting the existing "nom" (FE->FF) and
# "url_produit" (FF->FG) columns one position to the right.
$ws.Columns("FE").Insert()

# New header cell for the freshly inserted timestamp column
$ws.Range("FE1").Value = "2026-02-04 09:32:21"

# Populate the new FE column for each product row: carry forward the most
# recent known price (same as column FD) - except row 2 which records the
# newly observed price of 39.83. Rows with no price history (FD blank) are
# left blank as well.
$ws.Range("FE2").Value = 39.83
$ws.Range("FE3").Value = 169.95
$ws.Range("FE4").Value = 169.95
$ws.Range("FE5").Value = 199.95
$ws.Range("FE6").Value = 199.95
$ws.Range("FE7").Value = 199.95
$ws.Range("FE8").Value = 229.95
$ws.Range("FE9").Value = 249.95
$ws.Range("FE10").Value = 299.95
$ws.Range("FE11").Value = 619
$ws.Range("FE12").Value = 659
$ws.Range("FE13").Value = 659
$ws.Range("FE14").Value = 749
$ws.Range("FE15").Value = 809
$ws.Range("FE16").Value = 809
$ws.Range("FE17").Value = 809
$ws.Range("FE18").Value = 809
$ws.Range("FE19").Value = 809
$ws.Range("FE20").Value = 849
$ws.Range("FE21").Value = 899
$ws.Range("FE22").Value = 899
$ws.Range("FE23").Value = 909
$ws.Range("FE24").Value = 909
$ws.Range("FE25").Value = 909
$ws.Range("FE26").Value = 909
$ws.Range("FE27").Value = 969
$ws.Range("FE28").Value = 969
$ws.Range("FE29").Value = 969
$ws.Range("FE30").Value = 969
$ws.Range("FE31").Value = 969
$ws.Range("FE32").Value = 999
$ws.Range("FE33").Value = 999
$ws.Range("FE34").Value = 1039
$ws.Range("FE35").Value = 1039
$ws.Range("FE36").Value = 1079
$ws.Range("FE37").Value = 1079
$ws.Range("FE38").Value = 1079
$ws.Range("FE39").Value = 1079
$ws.Range("FE40").Value = 1099
$ws.Range("FE41").Value = 1099
$ws.Range("FE42").Value = 1199
$ws.Range("FE43").Value = 1219
$ws.Range("FE44").Value = 1219
$ws.Range("FE45").Value = 1219
$ws.Range("FE46").Value = 1219
$ws.Range("FE47").Value = 1219
$ws.Range("FE48").Value = 1229
$ws.Range("FE49").Value = 1229
$ws.Range("FE50").Value = 1249
$ws.Range("FE51").Value = 1329
$ws.Range("FE52").Value = 1329
$ws.Range("FE53").Value = 1329
$ws.Range("FE54").Value = 1329
$ws.Range("FE55").Value = 1329
$ws.Range("FE56").Value = 1329
$ws.Range("FE57").Value = 1329
$ws.Range("FE58").Value = 1349
$ws.Range("FE59").Value = 1479
$ws.Range("FE60").Value = 1479
$ws.Range("FE61").Value = 1479
$ws.Range("FE62").Value = 1579
$ws.Range("FE63").Value = 1579
$ws.Range("FE64").Value = 1579
$ws.Range("FE65").Value = 1579
$ws.Range("FE66").Value = 1579
$ws.Range("FE67").Value = 1579
$ws.Range("FE68").Value = 1579
$ws.Range("FE69").Value = 1729
$ws.Range("FE70").Value = 1729
$ws.Range("FE71").Value = 1729
$ws.Range("FE72").Value = 1829
$ws.Range("FE73").Value = 1829
$ws.Range("FE74").Value = 1829
$ws.Range("FE75").Value = 1979
$ws.Range("FE76").Value = 1979
$ws.Range("FE77").Value = 1979
$ws.Range("FE78").Value = 2479
$ws.Range("FE79").Value = 2479
$ws.Range("FE80").Value = 2479
